$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had 5 data rows (rows 2-6). Only the first two data rows
# (rows 2 and 3) are kept; rows 4, 5 and 6 are removed entirely.
$ws.Rows("4:6").Delete()

# Update the ResourceEmail (column D) for the two remaining data rows.
$ws.Range("D2").Value() = "dev@gmail.com"
$ws.Range("D3").Value() = "tester@gmail.com"

# Rebuild the hyperlinks collection: remove every stale mailto hyperlink
# (the ones that referenced the deleted rows as well as the old e-mail
# addresses) and re-create just the two that remain, pointing at the new
# addresses.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:tester@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:dev@gmail.com")

# Re-apply the built-in "Hyperlink" cell style so D2/D3 keep using the
# workbook's existing Hyperlink style entry instead of a newly minted one.
$ws.Range("D2").Style = "Hyperlink"
$ws.Range("D3").Style = "Hyperlink"

# Move the active selection to B3, matching the saved view state.
$ws.Range("B3").Select()
